# Update "想去人数" (want-to-go count) figures in column F for a handful of
# rows on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 324
    3 = 81
    5 = 4744
    6 = 374
    9 = 734
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
